$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 14 (pushes old rows 14-38 down to 16-40),
# mirroring the weekly data refresh described in the commit message.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# New row 14: Femacal de La Calera, Bruselas (repollito), 2022-06-08
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = "Femacal de La Calera"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44720
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 100112035
$ws.Range("G14").Value = "Bruselas (repollito)"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 85
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15529
$ws.Range("N14").Value = "$/malla 15 kilos"
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 1035
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = "Hortaliza"

# New row 15: Femacal de La Calera, Bruselas (repollito), 2022-06-09
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44721
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 100112035
$ws.Range("G15").Value = "Bruselas (repollito)"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 130
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = "$/malla 15 kilos"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 967
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = "Hortaliza"
